# Update the "Update automatically" date/time footer placeholder (the
# datetimeFigureOut field cached on the Slide Master and on every Slide
# Layout) from 3/16/2023 to 3/17/2023, reflecting the Mar 17 standup per
# the commit message. PowerPoint refreshes this cached field text on the
# master + every layout that carries a date placeholder (ppPlaceholderDate
# = 16) whenever the deck is touched on a later day.

$p = $ppt.ActivePresentation
$newDate = "3/17/2023"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $ph = $shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
